$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final")

# Header for the averages block (merged G4:H4)
$ws.Range("G4").Value = "Averages"
$ws.Range("G4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("G4:H4").Merge() | Out-Null

# Row 5: Year 1 overall average
$ws.Range("G5").Value = "Year 1"
$ws.Range("H5").Formula = '=SUMIF(Table1[Year],"1",Table1[Total])/COUNTIF(Table1[Year],"1")'
$ws.Range("H5").NumberFormat = "0.00%"

# Row 6: Year 1, Semester 1 average
$ws.Range("G6").Value = "Semester 1"
$ws.Range("H6").Formula = '=(SUMIFS(Table1[Total], Table1[Semester],"1",Table1[Year],"1"))/(COUNTIFS(Table1[Semester],"1",Table1[Year],"1"))'
$ws.Range("H6").NumberFormat = "0.00%"

# Row 7: Year 1, Semester 2 average
$ws.Range("G7").Value = "Semester 2"
$ws.Range("H7").Formula = '=(SUMIFS(Table1[Total], Table1[Semester],"2",Table1[Year],"1"))/(COUNTIFS(Table1[Semester],"2",Table1[Year],"1"))'
$ws.Range("H7").NumberFormat = "0.00%"

# Row 8: Year 2 overall average
$ws.Range("G8").Value = "Year 2"
$ws.Range("H8").Formula = '=SUMIF(Table1[Year],"2",Table1[Total])/COUNTIF(Table1[Year],"2")'
$ws.Range("H8").NumberFormat = "0.00%"

# Row 9: Year 2, Semester 1 average
$ws.Range("G9").Value = "Semester 1"
$ws.Range("H9").Formula = '=(SUMIFS(Table1[Total], Table1[Semester],"1",Table1[Year],"2"))/(COUNTIFS(Table1[Semester],"1",Table1[Year],"2"))'
$ws.Range("H9").NumberFormat = "0.00%"

# Row 10: Year 2, Semester 2 average
$ws.Range("G10").Value = "Semester 2"
$ws.Range("H10").Formula = '=(SUMIFS(Table1[Total],Table1[Semester],"2",Table1[Year],"2"))/(COUNTIFS(Table1[Semester],"2",Table1[Year],"2"))'
$ws.Range("H10").NumberFormat = "0.00%"

# Row 11: Year 3 overall average
$ws.Range("G11").Value = "Year 3"
$ws.Range("H11").Formula = '=IFERROR((SUMIF(Table1[Year],"3",Table1[Total])/COUNTIF(Table1[Year],"3")),0)'
$ws.Range("H11").NumberFormat = "0.00%"

# Row 12: Year 3, Semester 1 (placeholder, no data yet)
$ws.Range("G12").Value = "Semester 1"
$ws.Range("H12").NumberFormat = "0.00%"

# Row 13: Year 3, Semester 2 (placeholder, no data yet)
$ws.Range("G13").Value = "Semester 2"
$ws.Range("H13").NumberFormat = "0.00%"

$ws.Columns.Item(7).ColumnWidth = 12.140625

$wb.Save()
